$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3449
$ws1.Range("F3").Value = 27
$ws1.Range("F4").Value = 71
$ws1.Range("F5").Value = 1771
$ws1.Range("F6").Value = 100
$ws1.Range("F7").Value = 340

# Sheet "全部类型" (sheet4): update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3449
$ws4.Range("F3").Value = 27
$ws4.Range("F4").Value = 71
$ws4.Range("F5").Value = 1771
$ws4.Range("F6").Value = 100
$ws4.Range("F8").Value = 340
